# Applies the "Repartition" document update:
#  - month name in the title: سبتمبر -> أكتوبر
#  - numerous figures across the two summary tables
#  - the spelled-out total amount in Arabic words

$d = $word.ActiveDocument

# NOTE: the host object model re-resolves table references lazily, so we
# always fetch $d.Tables.Item($TableIndex) fresh right before touching a
# cell instead of caching table objects across interleaved table access.
function Set-CellText {
    param([int]$TableIndex, [int]$Row, [int]$Col, [string]$OldText, [string]$NewText)
    $cell = $d.Tables.Item($TableIndex).Cell($Row, $Col)
    $r = $cell.Range
    # Drop the trailing end-of-cell marker so we only touch the visible text.
    $r.End = $r.End - 1
    if ($r.Text -ne $OldText) {
        throw "Unexpected cell text at table $TableIndex row $Row col $Col : expected '$OldText' got '$($r.Text)'"
    }
    $r.Text = $NewText
}

# --- Table 1 ---
# عين تموشنت
Set-CellText 1 2 3 "888" "893"
Set-CellText 1 2 4 "8 880 000,00" "8 930 000,00"
Set-CellText 1 2 5 "10 530 000,00" "10 590 000,00"

# سيدي بن عدة
Set-CellText 1 3 3 "165" "166"
Set-CellText 1 3 4 "1 650 000,00" "1 660 000,00"

# المالح
Set-CellText 1 4 3 "201" "203"
Set-CellText 1 4 4 "2 010 000,00" "2 030 000,00"
Set-CellText 1 4 5 "4 900 000,00" "4 940 000,00"

# شعبة اللحم
Set-CellText 1 5 3 "152" "153"
Set-CellText 1 5 4 "1 520 000,00" "1 530 000,00"

# تارقة
Set-CellText 1 6 3 "89" "90"
Set-CellText 1 6 4 "890 000,00" "900 000,00"

# العامرية
Set-CellText 1 8 3 "194" "195"
Set-CellText 1 8 4 "1 940 000,00" "1 950 000,00"

# المساعيد
Set-CellText 1 12 3 "62" "61"
Set-CellText 1 12 4 "620 000,00" "610 000,00"

# حمام بوحجر
Set-CellText 1 13 3 "356" "359"
Set-CellText 1 13 4 "3 560 000,00" "3 590 000,00"
Set-CellText 1 13 5 "4 710 000,00" "4 750 000,00"

# وادي برقش
Set-CellText 1 16 3 "41" "42"
Set-CellText 1 16 4 "410 000,00" "420 000,00"

# عين الأربعاء
Set-CellText 1 17 5 "3 900 000,00" "3 910 000,00"

# وادي الصباح
Set-CellText 1 19 3 "105" "106"
Set-CellText 1 19 4 "1 050 000,00" "1 060 000,00"

# عين الكيحل
Set-CellText 1 21 3 "92" "94"
Set-CellText 1 21 4 "920 000,00" "940 000,00"
Set-CellText 1 21 5 "3 750 000,00" "3 770 000,00"

# المجموع (Table 1 total)
Set-CellText 1 25 3 "3271" "3288"
Set-CellText 1 25 4 "32 710 000,00" "32 880 000,00"
Set-CellText 1 25 5 "32 710 000,00" "32 880 000,00"

# --- Table 2 ---
# بني صاف
Set-CellText 2 2 3 "477" "475"
Set-CellText 2 2 4 "4 770 000,00" "4 750 000,00"

# سيدي الصافي
Set-CellText 2 3 3 "83" "84"
Set-CellText 2 3 4 "830 000,00" "840 000,00"

# الأمير عبد القادر
Set-CellText 2 4 3 "54" "55"
Set-CellText 2 4 4 "540 000,00" "550 000,00"

# المجموع العام (grand total)
Set-CellText 2 8 3 "4156" "4173"
Set-CellText 2 8 4 "41 560 000,00" "41 730 000,00"
Set-CellText 2 8 5 "41 560 000,00" "41 730 000,00"

# --- Non-table text ---
# Title month name
$d.Content.Find.Execute("سبتمبر", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "أكتوبر", 2) | Out-Null

# Spelled-out total amount in Arabic words
$d.Content.Find.Execute("واحد وأربعون مليون وخمسمئة وستون ألف", $false, $true, $false, $false, $false, `
                         $true, 1, $false, "واحد وأربعون مليون وسبعمئة وثلاثون ألف", 2) | Out-Null
